$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D ("Data de Nascimento"), shifting existing D..L to E..M
$ws.Columns("D:D").Insert()

# Header row
$ws.Cells.Item(1, 1).Value = "ID"
$ws.Cells.Item(1, 2).Value = "Nome"
$ws.Cells.Item(1, 3).Value = "Idade"
$ws.Cells.Item(1, 4).Value = "Data de Nascimento"
$ws.Cells.Item(1, 5).Value = "Gênero"
$ws.Cells.Item(1, 6).Value = "Etnia"
$ws.Cells.Item(1, 7).Value = "Nome da Mãe"
$ws.Cells.Item(1, 8).Value = "Educação"
$ws.Cells.Item(1, 9).Value = "Ocupação"
$ws.Cells.Item(1, 10).Value = "Telefone"
$ws.Cells.Item(1, 11).Value = "Celular"
$ws.Cells.Item(1, 12).Value = "CPF"
$ws.Cells.Item(1, 13).Value = "CEP"

# Data rows
# Row 2
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Andréa Thais Pires"
$ws.Cells.Item(2, 3).Value = 18
$ws.Cells.Item(2, 4).Value = "'27/10/2005"
$ws.Cells.Item(2, 5).Value = "Feminino"
$ws.Cells.Item(2, 6).Value = "Preta"
$ws.Cells.Item(2, 7).Value = "Márcia Thais"
$ws.Cells.Item(2, 8).Value = "Em idade escolar."
$ws.Cells.Item(2, 9).Value = "Empregado: Setor privado (sem CLT)"
$ws.Cells.Item(2, 10).Value = "73 5858-9702"
$ws.Cells.Item(2, 11).Value = "73 93599-6977"
$ws.Cells.Item(2, 12).Value = "819.931.870-81"
$ws.Cells.Item(2, 13).Value = "'757589557"

# Row 3
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "Hugo Lopes Gomes"
$ws.Cells.Item(3, 3).Value = 24
$ws.Cells.Item(3, 4).Value = "'21/01/1999"
$ws.Cells.Item(3, 5).Value = "Masculino"
$ws.Cells.Item(3, 6).Value = "Branco"
$ws.Cells.Item(3, 7).Value = "Patrícia Silva Lopes Gomes"
$ws.Cells.Item(3, 8).Value = "Em idade escolar."
$ws.Cells.Item(3, 9).Value = "Desocupado"
$ws.Cells.Item(3, 10).Value = "97 7172-6076"
$ws.Cells.Item(3, 11).Value = "97 95601-1925"
$ws.Cells.Item(3, 12).Value = "559.703.320-46"
$ws.Cells.Item(3, 13).Value = "'617568248"

# Row 4
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "Ana Ribeira Campos Pires Dias"
$ws.Cells.Item(4, 3).Value = 31
$ws.Cells.Item(4, 4).Value = "'30/07/1992"
$ws.Cells.Item(4, 5).Value = "Feminino"
$ws.Cells.Item(4, 6).Value = "Parda"
$ws.Cells.Item(4, 7).Value = "Cristina Pires"
$ws.Cells.Item(4, 8).Value = "Fundamental completo"
$ws.Cells.Item(4, 9).Value = "Autonomo: Sem CNPJ"
$ws.Cells.Item(4, 10).Value = "84 9500-9043"
$ws.Cells.Item(4, 11).Value = "84 90904-3479"
$ws.Cells.Item(4, 12).Value = "797.005.720-99"
$ws.Cells.Item(4, 13).Value = "'118920320"

# Row 5
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "Wagner Elias Gomes"
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = "'11/06/2022"
$ws.Cells.Item(5, 5).Value = "Masculino"
$ws.Cells.Item(5, 6).Value = "Branco"
$ws.Cells.Item(5, 7).Value = "Ana Gomes"
$ws.Cells.Item(5, 8).Value = "Em idade escolar."
$ws.Cells.Item(5, 9).Value = "Fora da força de trabalho"
$ws.Cells.Item(5, 10).Value = "64 8697-2935"
$ws.Cells.Item(5, 11).Value = "64 94734-0231"
$ws.Cells.Item(5, 12).Value = "330.173.150-37"
$ws.Cells.Item(5, 13).Value = "'857197361"

# Row 6
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "Vicente Ribeira Pinto"
$ws.Cells.Item(6, 3).Value = 49
$ws.Cells.Item(6, 4).Value = "'16/02/1974"
$ws.Cells.Item(6, 5).Value = "Masculino"
$ws.Cells.Item(6, 6).Value = "Branco"
$ws.Cells.Item(6, 7).Value = "Tatiane Leticia Pinto Ribeira"
$ws.Cells.Item(6, 8).Value = "Fundamental incompleto"
$ws.Cells.Item(6, 9).Value = "Autonomo: Com CNPJ"
$ws.Cells.Item(6, 10).Value = "31 3328-2454"
$ws.Cells.Item(6, 11).Value = "31 94760-8996"
$ws.Cells.Item(6, 12).Value = "635.867.450-10"
$ws.Cells.Item(6, 13).Value = "'081768997"

# Row 7
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "Fábio Ramos Oliveira"
$ws.Cells.Item(7, 3).Value = 23
$ws.Cells.Item(7, 4).Value = "'06/02/2000"
$ws.Cells.Item(7, 5).Value = "Masculino"
$ws.Cells.Item(7, 6).Value = "Pardo"
$ws.Cells.Item(7, 7).Value = "Isabela Maria Lopes Ramos Oliveira"
$ws.Cells.Item(7, 8).Value = "Em idade escolar."
$ws.Cells.Item(7, 9).Value = "Empregado: Setor ublico (estatutário ou militar)"
$ws.Cells.Item(7, 10).Value = "84 9966-3197"
$ws.Cells.Item(7, 11).Value = "84 90064-3826"
$ws.Cells.Item(7, 12).Value = "132.429.830-84"
$ws.Cells.Item(7, 13).Value = "'406990919"

# Row 8
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Karla Monteiro"
$ws.Cells.Item(8, 3).Value = 45
$ws.Cells.Item(8, 4).Value = "'19/05/1978"
$ws.Cells.Item(8, 5).Value = "Feminino"
$ws.Cells.Item(8, 6).Value = "Branca"
$ws.Cells.Item(8, 7).Value = "Patrícia Monteiro"
$ws.Cells.Item(8, 8).Value = "Superior completo"
$ws.Cells.Item(8, 9).Value = "Empregado: Trabalhador doméstico (sem CLT)"
$ws.Cells.Item(8, 10).Value = "94 3807-9995"
$ws.Cells.Item(8, 11).Value = "94 99156-7266"
$ws.Cells.Item(8, 12).Value = "709.487.240-21"
$ws.Cells.Item(8, 13).Value = "'454602676"

# Row 9
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Patrícia Costa"
$ws.Cells.Item(9, 3).Value = 42
$ws.Cells.Item(9, 4).Value = "'28/05/1981"
$ws.Cells.Item(9, 5).Value = "Feminino"
$ws.Cells.Item(9, 6).Value = "Parda"
$ws.Cells.Item(9, 7).Value = "Tânia Correia Costa"
$ws.Cells.Item(9, 8).Value = "Sem instrução"
$ws.Cells.Item(9, 9).Value = "Empregado: Setor privado (CLT)"
$ws.Cells.Item(9, 10).Value = "17 7852-1673"
$ws.Cells.Item(9, 11).Value = "17 96447-0527"
$ws.Cells.Item(9, 12).Value = "160.269.100-26"
$ws.Cells.Item(9, 13).Value = "'855671622"

# Row 10
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "Eduarda Oliveira"
$ws.Cells.Item(10, 3).Value = 22
$ws.Cells.Item(10, 4).Value = "'30/04/2001"
$ws.Cells.Item(10, 5).Value = "Feminino"
$ws.Cells.Item(10, 6).Value = "Parda"
$ws.Cells.Item(10, 7).Value = "Márcia Oliveira"
$ws.Cells.Item(10, 8).Value = "Em idade escolar."
$ws.Cells.Item(10, 9).Value = "Autonomo: Sem CNPJ"
$ws.Cells.Item(10, 10).Value = "75 3950-9309"
$ws.Cells.Item(10, 11).Value = "75 93548-4325"
$ws.Cells.Item(10, 12).Value = "077.496.700-50"
$ws.Cells.Item(10, 13).Value = "'161352608"

# Row 11
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Cátia Silva Dias"
$ws.Cells.Item(11, 3).Value = 54
$ws.Cells.Item(11, 4).Value = "'10/03/1969"
$ws.Cells.Item(11, 5).Value = "Feminino"
$ws.Cells.Item(11, 6).Value = "Preta"
$ws.Cells.Item(11, 7).Value = "Amanda Dias"
$ws.Cells.Item(11, 8).Value = "Fundamental incompleto"
$ws.Cells.Item(11, 9).Value = "Autonomo: Com CNPJ"
$ws.Cells.Item(11, 10).Value = "28 6071-3312"
$ws.Cells.Item(11, 11).Value = "28 91214-7612"
$ws.Cells.Item(11, 12).Value = "804.608.590-13"
$ws.Cells.Item(11, 13).Value = "'747872635"
